# Populate the (previously empty) sheet with a small 2x2 table:
#   A: a url text value, repeated on both rows
#   B: a date-looking text value - force text formatting first so
#      Excel doesn't reinterpret "2023-06-09" as a date serial number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1:B2").NumberFormat = "@"

$ws.Range("A1").Value = "www.google.com"
$ws.Range("B1").Value = "2023-06-09"
$ws.Range("A2").Value = "www.google.com"
$ws.Range("B2").Value = "2023-06-09"
